$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the shared string label for B6 (row 6, "Additional lipid extraction efficiency [%]" -> "Bagasse lipid extraction efficiency [%]")
$ws.Range("B6").Value = "Bagasse lipid extraction efficiency [%]"

# Update Spearman correlation values (rows 4-13, columns C,D,E,G,H,I,J)
$ws.Range("C4").Value = 0.1549868746718668
$ws.Range("D4").Value = 0.9549208730218257
$ws.Range("E4").Value = -0.7605295132378309
$ws.Range("G4").Value = -0.04782869571739294
$ws.Range("H4").Value = 0.4709592739818496
$ws.Range("I4").Value = 0.4501942548563714
$ws.Range("J4").Value = 0.8443291082277057
$ws.Range("C5").Value = 0.04119402985074627
$ws.Range("D5").Value = 0.09569489237230933
$ws.Range("E5").Value = 0.03064726618165454
$ws.Range("G5").Value = 0.09407185179629492
$ws.Range("H5").Value = 0.1516552913822846
$ws.Range("I5").Value = 0.1352568814220356
$ws.Range("J5").Value = 0.2391569789244731
$ws.Range("C6").Value = 0.07290332258306458
$ws.Range("D6").Value = 0.1855126378159454
$ws.Range("E6").Value = -0.02396009900247507
$ws.Range("G6").Value = 0.1831845796144904
$ws.Range("H6").Value = 0.1491832295807395
$ws.Range("I6").Value = 0.0005745143628590716
$ws.Range("J6").Value = 0.04837170929273232
$ws.Range("C7").Value = 0.0983844596114903
$ws.Range("D7").Value = 0.08411160279006977
$ws.Range("E7").Value = 0.5457226430660768
$ws.Range("G7").Value = 0.3583769594239857
$ws.Range("H7").Value = 0.6976209405235132
$ws.Range("I7").Value = 0.8364224105602641
$ws.Range("J7").Value = -0.08258756468911724
$ws.Range("C8").Value = 0.8228275706892674
$ws.Range("D8").Value = 0.04603165079126979
$ws.Range("E8").Value = 0.001239030975774395
$ws.Range("G8").Value = 0.01265431635790895
$ws.Range("H8").Value = 0.07492387309682742
$ws.Range("I8").Value = 0.06387759693992351
$ws.Range("J8").Value = 0.01538888472211805
$ws.Range("C9").Value = 0.3799579989499738
$ws.Range("D9").Value = -0.09046426160654017
$ws.Range("E9").Value = 0.03308932723318084
$ws.Range("G9").Value = -0.03423085577139429
$ws.Range("H9").Value = -0.06160954023850598
$ws.Range("I9").Value = -0.02327308182704568
$ws.Range("J9").Value = -0.06547513687842198
$ws.Range("C10").Value = 0.03581639540988525
$ws.Range("D10").Value = 0.01844596114902873
$ws.Range("E10").Value = -0.03937298432460812
$ws.Range("G10").Value = 0.01673891847296183
$ws.Range("H10").Value = -0.005610140253506338
$ws.Range("I10").Value = -0.01525238130953274
$ws.Range("J10").Value = 0.003933098327458187
$ws.Range("C11").Value = -0.01971649291232281
$ws.Range("D11").Value = 0.08459011475286883
$ws.Range("E11").Value = 0.03330683267081678
$ws.Range("G11").Value = 0.05852096302407562
$ws.Range("H11").Value = 0.1238580964524113
$ws.Range("I11").Value = 0.1296977424435611
$ws.Range("J11").Value = 0.005083627090677268
$ws.Range("C12").Value = 0.1048961224030601
$ws.Range("D12").Value = 0.1101987549688742
$ws.Range("E12").Value = 0.3554188854721369
$ws.Range("G12").Value = 0.8372519312982827
$ws.Range("H12").Value = 0.3183154578864472
$ws.Range("I12").Value = -0.08373959348983726
$ws.Range("J12").Value = 0.01337883447086177
$ws.Range("C13").Value = -0.2368589214730369
$ws.Range("D13").Value = 0.003964599114977875
$ws.Range("E13").Value = -0.02074551863796595
$ws.Range("G13").Value = -0.03131178279456987
$ws.Range("H13").Value = 0.03660091502287557
$ws.Range("I13").Value = 0.05649141228530714
$ws.Range("J13").Value = 0.000006000150003750094
